{"js": "// Update Alex Morgan resume content: Professional Summary, Education date,\n// Skills list, and Experience section (title + 3 bullet points) to reflect\n// the \"Tech Solutions\" role and refreshed accomplishments.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map each paragraph's current (old) text to its replacement (new) text so\n// we can locate the right paragraph robustly regardless of index drift.\nconst replacements = [\n  {\n    oldText:\n      \"Detail-oriented Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in predictive analytics and data modeling. Highly skilled in Python, SQL, and machine learning, with a proven track record of enhancing decision-making efficiency by 15% through advanced data analysis techniques. Demonstrates strong problem-solving capabilities and a knack for translating complex data insights into actionable business strategies. Certified in data science, with extensive experience in quantitative analysis and statistical modeling. Eager to leverage technical expertise and analytical skills in a dynamic new role to drive business success.\",\n    newText:\n      \"Accomplished Data Scientist with a Bachelor of Science in Computer Science from Stanford University, bringing a robust background in machine learning, predictive analytics, and generating business insights. Demonstrated success at Tech Solutions, where I enhanced forecast accuracy by 20% through sophisticated modeling techniques and improved operational efficiency by 15% via strategic data analysis. Proficient in Python, R, SQL, and adept in utilizing machine learning frameworks like TensorFlow and Scikit-Learn. Certified in data science with extensive experience in transforming complex datasets into actionable insights to propel business objectives. Seeking to apply my analytical prowess and technical skills in a dynamic new setting to drive further business successes.\",\n  },\n  {\n    oldText: \"B.S. in Computer Science, Stanford University, 2022\",\n    newText: \"B.S. in Computer Science, Stanford University, 2020\",\n  },\n  {\n    oldText:\n      \"Python, SQL, Machine Learning, Predictive Analytics, Data Visualization, Statistical Modeling, Quantitative Analysis, Problem Solving\",\n    newText:\n      \"Python, R, SQL, Machine Learning, Predictive Analytics, TensorFlow, Scikit-Learn, Data Visualization, Tableau\",\n  },\n  {\n    oldText: \"Data Analyst at TechCorp Analytics (Jan 2022 - Present)\",\n    newText: \"Data Analyst at Tech Solutions (Jan 2021 \\u2013 Feb 2023)\",\n  },\n  {\n    oldText:\n      \"Developed and implemented predictive models using Python to analyze large-scale datasets, delivering actionable insights that significantly impacted business strategies.\",\n    newText:\n      \"Developed and refined advanced forecasting models using Python and TensorFlow, which increased sales forecast accuracy by 20%.\",\n  },\n  {\n    oldText:\n      \"Enhanced data interpretation capabilities leading to a 15% increase in decision-making efficiency for key business units, through the use of advanced statistical modeling and machine learning techniques.\",\n    newText:\n      \"Designed and implemented interactive dashboards and visualizations in Tableau, significantly enhancing management's decision-making capabilities.\",\n  },\n  {\n    oldText:\n      \"Collaborated with cross-functional teams to refine data collection and analysis processes, improving data reliability by 20% and ensuring data integrity in project outcomes.\",\n    newText:\n      \"Performed in-depth data analysis to extract strategic insights, thereby boosting operational efficiency by 15% and supporting key business initiatives.\",\n  },\n];\n\nfor (const item of paragraphs.items) {\n  const match = replacements.find((r) => item.text === r.oldText);\n  if (match) {\n    item.insertText(match.newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update Alex Morgan resume content: Professional Summary, Education date,\n# Skills list, and Experience section (title + 3 bullet points) to reflect\n# the \"Tech Solutions\" role and refreshed accomplishments.\n\n$d = $word.ActiveDocument\n$enDash = [char]0x2013\n\n$replacements = @(\n    @{\n        Old = \"Detail-oriented Data Scientist with a Bachelor of Science in Computer Science from Stanford University, specializing in predictive analytics and data modeling. Highly skilled in Python, SQL, and machine learning, with a proven track record of enhancing decision-making efficiency by 15% through advanced data analysis techniques. Demonstrates strong problem-solving capabilities and a knack for translating complex data insights into actionable business strategies. Certified in data science, with extensive experience in quantitative analysis and statistical modeling. Eager to leverage technical expertise and analytical skills in a dynamic new role to drive business success.\"\n        New = \"Accomplished Data Scientist with a Bachelor of Science in Computer Science from Stanford University, bringing a robust background in machine learning, predictive analytics, and generating business insights. Demonstrated success at Tech Solutions, where I enhanced forecast accuracy by 20% through sophisticated modeling techniques and improved operational efficiency by 15% via strategic data analysis. Proficient in Python, R, SQL, and adept in utilizing machine learning frameworks like TensorFlow and Scikit-Learn. Certified in data science with extensive experience in transforming complex datasets into actionable insights to propel business objectives. Seeking to apply my analytical prowess and technical skills in a dynamic new setting to drive further business successes.\"\n    },\n    @{\n        Old = \"B.S. in Computer Science, Stanford University, 2022\"\n        New = \"B.S. in Computer Science, Stanford University, 2020\"\n    },\n    @{\n        Old = \"Python, SQL, Machine Learning, Predictive Analytics, Data Visualization, Statistical Modeling, Quantitative Analysis, Problem Solving\"\n        New = \"Python, R, SQL, Machine Learning, Predictive Analytics, TensorFlow, Scikit-Learn, Data Visualization, Tableau\"\n    },\n    @{\n        Old = \"Data Analyst at TechCorp Analytics (Jan 2022 - Present)\"\n        New = \"Data Analyst at Tech Solutions (Jan 2021 \" + $enDash + \" Feb 2023)\"\n    },\n    @{\n        Old = \"Developed and implemented predictive models using Python to analyze large-scale datasets, delivering actionable insights that significantly impacted business strategies.\"\n        New = \"Developed and refined advanced forecasting models using Python and TensorFlow, which increased sales forecast accuracy by 20%.\"\n    },\n    @{\n        Old = \"Enhanced data interpretation capabilities leading to a 15% increase in decision-making efficiency for key business units, through the use of advanced statistical modeling and machine learning techniques.\"\n        New = \"Designed and implemented interactive dashboards and visualizations in Tableau, significantly enhancing management's decision-making capabilities.\"\n    },\n    @{\n        Old = \"Collaborated with cross-functional teams to refine data collection and analysis processes, improving data reliability by 20% and ensuring data integrity in project outcomes.\"\n        New = \"Performed in-depth data analysis to extract strategic insights, thereby boosting operational efficiency by 15% and supporting key business initiatives.\"\n    }\n)\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\")\n    foreach ($rep in $replacements) {\n        if ($text -eq $rep.Old) {\n            $p.Range.Text = $rep.New\n            break\n        }\n    }\n}\n"}
